# Finished Week 13 logging
# Update Home-row (row 2) Short/Deep attempt/completion counts
# on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 179
$wsOff.Range("C2").Value = 133
$wsOff.Range("D2").Value = 53
$wsOff.Range("E2").Value = 25

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 178
$wsDef.Range("C2").Value = 115
$wsDef.Range("D2").Value = 32
$wsDef.Range("E2").Value = 14
